$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.878.17"
$ws.Range("E2").Value = "  -2.00%  "

# Row 3
$ws.Range("D3").Value = "1.632.84"
$ws.Range("E3").Value = "  -2.17%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.41%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "

# Row 6
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.32%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.012"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.42%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2564"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.36%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06408"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.01%  "

# Row 12
$ws.Range("D12").Value = "1.640.98"
$ws.Range("E12").Value = "  -1.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.244"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.48%  "

# Row 14
$ws.Range("D14").Value = "1.861.30"
$ws.Range("E14").Value = "  -1.95%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5435"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.07%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7916"
$ws.Range("E16").Value = "  -1.77%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.10%  "

# Row 18
$ws.Range("D18").Value = "25.906.79"
$ws.Range("E18").Value = "  -2.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.68%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.296"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.93%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.987"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.942"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.69%  "

# Row 24
$ws.Range("E24").Value = "  +0.44%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.970"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.97%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1146"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.01%  "

# Row 28
$ws.Range("E28").Value = "  -0.66%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.775"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05068"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.75%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.241"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.257"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.542"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.22%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.345"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8893"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.605"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.74%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5638"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.70%  "

# Row 39
$ws.Range("D39").Value = "1.133.46"
$ws.Range("E39").Value = "  -1.90%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01561"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.85%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.580"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.011"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.42%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.635"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8166"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "

# Row 46
$ws.Range("D46").Value = "1.772.21"
$ws.Range("E46").Value = "  -1.93%  "

# Row 47
$ws.Range("E47").Value = "  +0.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4532"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.012"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.44%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05023"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.74%  "
